# working full input parsing
#
# Adds the first "guest" row of parsed data to the guests sheet (mirrors
# the row already present on the players sheet) and updates the row-2
# selection on the players sheet to span the full row.

$wb = $excel.ActiveWorkbook

# --- guests sheet: populate row 2 with the parsed guest entry ---
$guests = $wb.Worksheets.Item("guests")
$guests.Range("A2").Value = "Botond"
$guests.Range("B2").Value = "Yedo"
$guests.Range("C2").Value = "Hegemony"
$guests.Range("D2").Value = "Ark Nova"
$guests.Range("E2").Value = "Peter"
$guests.Range("F2").Value = "Yes - to play something different (pls specify in comments section)"
$guests.Range("I2").Value = "Item Type"
$guests.Range("J2").Value = "company/"

# Column widths widen to fit the newly entered content (mirrors Excel's
# "best fit" recompute that happens automatically once longer values are
# entered in a bestFit column).
$guests.Columns.Item(1).ColumnWidth = 6.666666666666667
$guests.Columns.Item(6).ColumnWidth = 60.666666666666664
$guests.Columns.Item(10).ColumnWidth = 9

# --- players sheet: select the whole of row 2 (A2:XFD2) ---
$players = $wb.Worksheets.Item("players")
$players.Rows.Item(2).Select() | Out-Null

# restore guests as the active/visible sheet (matches the saved workbook state)
$guests.Activate()
